$wb = $excel.ActiveWorkbook

# =====================================================================
# Step 1: rename the existing "timeline" sheet to "timeline raw data"
# =====================================================================
$tl = $wb.Worksheets.Item("timeline")
$tl.Name = "timeline raw data"

# =====================================================================
# Step 2: add summary CONCATENATE formulas in column F of each block's
# header row (one formula per job/school/cause entry).
# =====================================================================
$tl.Range("F12").Formula = "=CONCATENATE(C12,E12,D12,C13,E13,D13,C14,E14,D14,C15,E15,D15,C16,E16,D16)"
$tl.Range("F17").Formula = "=CONCATENATE(C17,E17,D17,C18,E18,D18,C19,E19,D19,C20,E20,D20,C21,E21,D21)"
$tl.Range("F22").Formula = "=CONCATENATE(C22,E22,D22,C23,E23,D23,C24,E24,D24,C25,E25,D25,C26,E26,D26,C27,E27,D27)"
$tl.Range("F28").Formula = "=CONCATENATE(C28,E28,D28,C29,E29,D29,C30,E30,D30,C31,E31,D31,C32,E32,D32,C33,E33,D33)"
$tl.Range("F34").Formula = "=CONCATENATE(C34,E34,D34,C35,E35,D35,C36,E36,D36,C37,E37,D37,C38,E38,D38,C39,E39,D39)"
$tl.Range("F40").Formula = "=CONCATENATE(C40,E40,D40,C41,E41,D41,C42,E42,D42,C43,E43,D43,C44,E44,D44,C45,E45,D45)"
$tl.Range("F46").Formula = "=CONCATENATE(C46,E46,D46,C47,E47,D47,C48,E48,D48,C49,E49,D49,C50,E50,D50,C51,E51,D51)"
$tl.Range("F52").Formula = "=CONCATENATE(C52,E52,D52,C53,E53,D53,C54,E54,D54,C55,E55,D55,C56,E56,D56)"
$tl.Range("F57").Formula = "=CONCATENATE(C57,E57,D57,C58,E58,D58,C59,E59,D59,C60,E60,D60,C61,E61,D61)"
$tl.Range("F62").Formula = "=CONCATENATE(C62,E62,D62,C63,E63,D63,C64,E64,D64,C65,E65,D65,C66,E66,D66,C67,E67,D67)"
$tl.Range("F68").Formula = "=CONCATENATE(C68,E68,D68,C69,E69,D69,C70,E70,D70,C71,E71,D71,C72,E72,D72)"

# set the selection on the raw-data (old) sheet to match target
$tl.Range("B20:C20").Select()

# =====================================================================
# Step 3: insert a brand new "Sheet1" before it (becomes the new first
# tab) -- this holds the clean, pasted-in raw timeline data.
# =====================================================================
$raw = $wb.Worksheets.Add()
$raw.Name = "Sheet1"

$raw.Columns.Item(1).ColumnWidth = 26.42578125
$raw.Columns.Item(2).ColumnWidth = 9.28515625
$raw.Columns.Item(3).ColumnWidth = 12.140625

# Header row
$raw.Range("A1").Value = "cat"
$raw.Range("B1").Value = "beg"
$raw.Range("C1").Value = "end"
$raw.Range("D1").Value = "des"

# Pre-format the date columns with the "d-mmm-yy" built-in format
# (must happen before assigning .Value so the date style sticks cleanly)
$raw.Range("B2:C23").NumberFormat = "d-mmm-yy"

# Row 2: undergrad
$raw.Range("A2").Value = "undergrad"
$raw.Range("B2").Value = (Get-Date -Year 2005 -Month 9 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C2").Value = (Get-Date -Year 2009 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D2").Value = "<h3>uc davis</h3><h4>economics major</h4><h5>sep. 05 to mar. 09</h5><p>2009: b.a. in economics</p><p class='placeLabel'>DAVIS, CA</p>"

# Row 3: growth marketing - catchafire
$raw.Range("A3").Value = "growth marketing"
$raw.Range("B3").Value = (Get-Date -Year 2014 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C3").Value = (Get-Date -Year 2015 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D3").Value = "<h3>catchafire</h3><h4>search marketing strategy</h4><h5>dec. 14 to present</h5><p>provide customer acquisition and analytical insights for social causes and change makers</p><p class='placeLabel'>SAN FRANCISCO BAY AREA</p>"

# Row 4: growth marketing - inflection
$raw.Range("A4").Value = "growth marketing"
$raw.Range("B4").Value = (Get-Date -Year 2013 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C4").Value = (Get-Date -Year 2014 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D4").Value = "<h3>inflection</h3><h4>senior search engine marketing specialist</h4><h5>aug. 13 to dec. 14</h5><p>acquired high converting traffic via paid acquisition channels</p><p>optimized campaigns based on customer lifetime value</p><p class='placeLabel'>SAN FRANCISCO BAY AREA</p>"

# Row 5: growth marketing - brighter collective
$raw.Range("A5").Value = "growth marketing"
$raw.Range("B5").Value = (Get-Date -Year 2012 -Month 11 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C5").Value = (Get-Date -Year 2013 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D5").Value = "<h3>brighter collective</h3><h4>media analyst</h4><h5>nov. 12 to aug. 13</h5><p>led customer acquisition search and display campaigns</p><p>optimized text ad copy, keywords, and built campaigns by hand</p><p class='placeLabel'>LOS ANGELES</p>"

# Row 6: growth marketing - ticketmaster
$raw.Range("A6").Value = "growth marketing"
$raw.Range("B6").Value = (Get-Date -Year 2012 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C6").Value = (Get-Date -Year 2012 -Month 11 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D6").Value = "<h3>ticketmaster</h3><h4>search analyst</h4><h5>mar. 12 to nov. 12</h5><p>analyzed key performance indicators for search marketing, search engine optimization, and social</p><p>toolbelt included adwords, sitecatalyst, and webmaster tools</p><p class='placeLabel'>HOLLYWOOD</p>"

# Row 7: growth marketing - pricegrabber (search marketing analyst)
$raw.Range("A7").Value = "growth marketing"
$raw.Range("B7").Value = (Get-Date -Year 2011 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C7").Value = (Get-Date -Year 2012 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D7").Value = "<h3>pricegrabber</h3><h4>search marketing analyst</h4><h5>jul. 11 to mar. 12</h5><p>performed reporting and campaign optimizations</p><p>communicated daily performance results to ceo and key people</p><p class='placeLabel'>LOS ANGELES</p>"

# Row 8: growth marketing - pricegrabber (client services representative)
$raw.Range("A8").Value = "growth marketing"
$raw.Range("B8").Value = (Get-Date -Year 2010 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C8").Value = (Get-Date -Year 2011 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D8").Value = "<h3>pricegrabber</h3><h4>client services representative</h4><h5>apr. 10 to jul. 11</h5><p>helped clients improve ROI with pay-per-click bid suggestions</p><p>upped merchant participation in value-add features such as conversion tracking</p><p class='placeLabel'>LOS ANGELES</p>"

# Row 9: giving - urban light
$raw.Range("A9").Value = "giving"
$raw.Range("B9").Value = (Get-Date -Year 2010 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C9").Value = (Get-Date -Year 2015 -Month 12 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D9").Value = "<h3>urban light</h3><h4>fundraising</h4><h5>apr. 2010 to present</h5><p>raise funds for a game changing organization</p><p class='placeLabel'>CHIANG MAI</p>"

# Row 10: awareness events - human trafficking awareness orgs
$raw.Range("A10").Value = "awareness events"
$raw.Range("B10").Value = (Get-Date -Year 2009 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C10").Value = (Get-Date -Year 2010 -Month 4 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D10").Value = "<h3>human trafficking awareness orgs</h3><h4>awareness event planning</h4><h5>aug. 09 to apr. 2010</h5><p>planned events, screened documentaries in public venues</p><p class='placeLabel'>DAVIS, CA</p>"

# Row 11: overseas - the sold project
$raw.Range("A11").Value = "overseas"
$raw.Range("B11").Value = (Get-Date -Year 2009 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C11").Value = (Get-Date -Year 2009 -Month 8 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D11").Value = "<h3>the sold project</h3><h4>overseas volunteer</h4><h5>jul. 09 to aug. 09</h5><p>taught english and did supply drops at local ngos</p><p>crossed the thai - burmese border gate</p><p class='placeLabel'>NORTHERN THAILAND</p>"

# Row 12: awareness events - various orgs
$raw.Range("A12").Value = "awareness events"
$raw.Range("B12").Value = (Get-Date -Year 2009 -Month 3 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("C12").Value = (Get-Date -Year 2009 -Month 7 -Day 1 -Hour 0 -Minute 0 -Second 0).Date
$raw.Range("D12").Value = "<h3>various orgs</h3><h4>awareness event planning</h4><h5>mar. 09 to jul. 09</h5><p>worked with with a great team and planned human trafficking awareness events</p><p class='placeLabel'>DAVIS, CA</p>"

# Rows 15-23: pre-formatted (date-styled) blank placeholder cells in B:C
$raw.Range("B15:C23").NumberFormat = "d-mmm-yy"

$raw.Range("A14").Select()
